# Apply cryptos list update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "61.528.66"
Set-TextValue $ws.Range("E2") "  +0.98%  "

Set-TextValue $ws.Range("D3") "3.452.55"
Set-TextValue $ws.Range("E3") "  +2.07%  "

Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  -0.05%  "

Set-TextValue $ws.Range("D5") "580.06"
Set-TextValue $ws.Range("E5") "  +1.34%  "

Set-TextValue $ws.Range("D6") "148.36"
Set-TextValue $ws.Range("E6") "  +8.55%  "

Set-TextValue $ws.Range("D7") "3.453.46"
Set-TextValue $ws.Range("E7") "  +2.15%  "

Set-TextValue $ws.Range("E8") "  +0.05%  "

Set-TextValue $ws.Range("E9") "  +0.56%  "

Set-TextValue $ws.Range("D10") "7.81"
Set-TextValue $ws.Range("E10") "  +3.70%  "

Set-TextValue $ws.Range("D11") "0.125"
Set-TextValue $ws.Range("E11") "  +0.71%  "

Set-TextValue $ws.Range("D12") "0.392"
Set-TextValue $ws.Range("E12") "  +1.04%  "

Set-TextValue $ws.Range("D13") "4.043.44"
Set-TextValue $ws.Range("E13") "  +2.18%  "

Set-TextValue $ws.Range("D14") "28.05"
Set-TextValue $ws.Range("E14") "  +6.17%  "

Set-TextValue $ws.Range("E15") "  -0.60%  "

Set-TextValue $ws.Range("E16") "  +1.52%  "

Set-TextValue $ws.Range("D17") "3.438.39"
Set-TextValue $ws.Range("E17") "  +1.68%  "

Set-TextValue $ws.Range("D18") "61.622.30"
Set-TextValue $ws.Range("E18") "  +0.86%  "

Set-TextValue $ws.Range("D19") "6.33"
Set-TextValue $ws.Range("E19") "  +8.45%  "

Set-TextValue $ws.Range("D20") "14.36"
Set-TextValue $ws.Range("E20") "  +2.45%  "

Set-TextValue $ws.Range("D21") "9.47"
Set-TextValue $ws.Range("E21") "  +0.98%  "

Set-TextValue $ws.Range("D22") "386.56"
Set-TextValue $ws.Range("E22") "  +2.55%  "

Set-TextValue $ws.Range("E23") "  +2.30%  "

Set-TextValue $ws.Range("D24") "3.595.13"
Set-TextValue $ws.Range("E24") "  +2.48%  "

Set-TextValue $ws.Range("D25") "72.73"
Set-TextValue $ws.Range("E25") "  +2.04%  "

Set-TextValue $ws.Range("E26") "  -0.05%  "

Set-TextValue $ws.Range("D27") "5.78"
Set-TextValue $ws.Range("E27") "  +0.79%  "

Set-TextValue $ws.Range("E28") "  -2.34%  "

Set-TextValue $ws.Range("E29") "  +7.60%  "

Set-TextValue $ws.Range("D30") "7.84"
Set-TextValue $ws.Range("E30") "  +3.92%  "

Set-TextValue $ws.Range("E31") "  -0.04%  "

Set-TextValue $ws.Range("E32") "  -13.85%  "

Set-TextValue $ws.Range("E33") "  +1.36%  "

Set-TextValue $ws.Range("E34") "  +1.04%  "

Set-TextValue $ws.Range("E35") "  +0.01%  "

Set-TextValue $ws.Range("D36") "23.99"
Set-TextValue $ws.Range("E36") "  +0.76%  "

Set-TextValue $ws.Range("D37") "7.09"
Set-TextValue $ws.Range("E37") "  +3.89%  "

Set-TextValue $ws.Range("E38") "  +0.40%  "

Set-TextValue $ws.Range("D39") "1.57"
Set-TextValue $ws.Range("E39") "  +2.45%  "

Set-TextValue $ws.Range("D40") "166.37"
Set-TextValue $ws.Range("E40") "  +0.88%  "

Set-TextValue $ws.Range("D41") "0.0791"
Set-TextValue $ws.Range("E41") "  +4.62%  "

Set-TextValue $ws.Range("D44") "4.53"
Set-TextValue $ws.Range("E44") "  +2.65%  "

Set-TextValue $ws.Range("E45") "  -0.03%  "

Set-TextValue $ws.Range("D46") "42.35"
Set-TextValue $ws.Range("E46") "  +1.84%  "

Set-TextValue $ws.Range("E47") "  +1.32%  "

Set-TextValue $ws.Range("D48") "2.609.43"
Set-TextValue $ws.Range("E48") "  +9.82%  "

Set-TextValue $ws.Range("E49") "  -3.84%  "

Set-TextValue $ws.Range("D50") "7.00"
Set-TextValue $ws.Range("E50") "  +2.90%  "

Set-TextValue $ws.Range("D51") "23.26"
Set-TextValue $ws.Range("E51") "  -0.58%  "

# Row 42/43: EnergySwap and Mantle swapped positions with new data
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D42") "0.797"
Set-TextValue $ws.Range("E42") "  +3.24%  "

$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D43") "26.00"
Set-TextValue $ws.Range("E43") "  +8.55%  "
